$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @(2, "61.482.35", "  +1.15%  "),
  @(3, "3.447.22", "  +1.63%  "),
  @(4, "0.999", "  -0.16%  "),
  @(5, "579.86", "  +0.82%  "),
  @(6, "148.83", "  +8.56%  "),
  @(7, "3.447.57", "  +1.61%  "),
  @(8, $null, "  +0.08%  "),
  @(9, $null, "  +0.63%  "),
  @(10, $null, "  +3.04%  "),
  @(11, "0.125", "  +1.04%  "),
  @(12, "0.392", "  +0.26%  "),
  @(13, "4.036.38", "  +1.21%  "),
  @(14, "28.11", "  +5.89%  "),
  @(15, $null, "  -0.43%  "),
  @(16, $null, "  +0.92%  "),
  @(17, "3.451.61", "  +1.22%  "),
  @(18, "61.580.68", "  +0.95%  "),
  @(19, "6.33", "  +7.94%  "),
  @(20, "14.40", "  +2.48%  "),
  @(21, "9.44", "  -0.31%  "),
  @(22, "387.33", "  +2.28%  "),
  @(23, "0.571", "  +2.18%  "),
  @(24, "3.591.66", "  +1.86%  "),
  @(25, "72.74", "  +2.19%  "),
  @(26, $null, "  -0.19%  "),
  @(27, "5.78", "  +0.65%  "),
  @(28, "0.0000123", "  -2.02%  "),
  @(29, "0.181", "  +7.27%  "),
  @(30, "7.81", "  +2.47%  "),
  @(31, $null, "  -0.06%  "),
  @(32, "1.54", "  -14.10%  "),
  @(33, $null, "  +0.37%  "),
  @(34, "2.18", "  +0.42%  "),
  @(35, $null, "  +0.01%  "),
  @(36, "23.98", "  +0.59%  "),
  @(37, "7.09", "  +2.59%  "),
  @(38, "5.26", "  +0.64%  "),
  @(39, $null, "  +1.55%  "),
  @(40, "165.92", "  +0.96%  "),
  @(41, "0.0794", "  +4.27%  "),
  @(42, "26.07", "  +6.69%  "),
  @(43, "0.795", "  +2.68%  "),
  @(44, $null, "  -0.24%  "),
  @(45, "4.50", "  +1.72%  "),
  @(46, "42.30", "  +1.78%  "),
  @(47, $null, "  +1.09%  "),
  @(48, "2.613.75", "  +7.77%  "),
  @(49, $null, "  -3.58%  "),
  @(50, "7.03", "  +3.08%  "),
  @(51, "23.10", "  -2.11%  ")
)

foreach ($item in $updates) {
    $row = $item[0]
    $dVal = $item[1]
    $eVal = $item[2]
    if ($dVal) {
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
    }
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $eVal
}

Write-Output "Done updating cryptos list"